# Updates cryptos list figures (price + 1h volume change) scraped on
# Wed Aug 14 13:42:03 UTC 2024. Also fixes the WrappedEther/ShibaInu row
# ordering (rows 17-18) and swaps the InjectiveProtocol row (51) for VeChain.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "60.532.08"
$c.ClearFormats()
$ws.Range("E2").Value = "  +2.69%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.695.97"
$c.ClearFormats()
$ws.Range("E3").Value = "  +2.14%  "

$ws.Range("E4").Value = "  +0.16%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "526.42"
$c.ClearFormats()
$ws.Range("E5").Value = "  +1.19%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "149.35"
$c.ClearFormats()
$ws.Range("E6").Value = "  +1.88%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +1.27%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.716.39"
$c.ClearFormats()
$ws.Range("E9").Value = "  +2.58%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.00"
$c.ClearFormats()
$ws.Range("E10").Value = "  +10.85%  "

$ws.Range("E11").Value = "  +0.24%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.340"
$c.ClearFormats()
$ws.Range("E12").Value = "  +1.39%  "

$ws.Range("E13").Value = "  +2.39%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.175.24"
$c.ClearFormats()
$ws.Range("E14").Value = "  +2.47%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "60.526.72"
$c.ClearFormats()
$ws.Range("E15").Value = "  +2.69%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "21.50"
$c.ClearFormats()
$ws.Range("E16").Value = "  +2.83%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.705.30"
$c.ClearFormats()
$ws.Range("E17").Value = "  +2.49%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.0000138"
$c.ClearFormats()
$ws.Range("E18").Value = "  +0.58%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "345.11"
$c.ClearFormats()
$ws.Range("E19").Value = "  -1.31%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.50"
$c.ClearFormats()
$ws.Range("E20").Value = "  +0.55%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.54"
$c.ClearFormats()
$ws.Range("E21").Value = "  +2.12%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.40"
$c.ClearFormats()
$ws.Range("E22").Value = "  +3.51%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.01%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "63.71"
$c.ClearFormats()
$ws.Range("E24").Value = "  +2.96%  "

$ws.Range("E25").Value = "  +4.10%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.418"
$c.ClearFormats()
$ws.Range("E26").Value = "  +0.51%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.ClearFormats()
$ws.Range("E27").Value = "  +0.17%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.0₃0820"
$c.ClearFormats()
$ws.Range("E28").Value = "  +1.74%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.29"
$c.ClearFormats()
$ws.Range("E29").Value = "  +3.01%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.81"
$c.ClearFormats()
$ws.Range("E30").Value = "  +7.88%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.ClearFormats()
$ws.Range("E31").Value = "  -0.05%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.60"
$c.ClearFormats()
$ws.Range("E32").Value = "  +1.38%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "19.06"
$c.ClearFormats()
$ws.Range("E33").Value = "  +0.80%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "150.54"
$c.ClearFormats()
$ws.Range("E34").Value = "  +0.77%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.26"
$c.ClearFormats()
$ws.Range("E35").Value = "  +5.62%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.22"
$c.ClearFormats()
$ws.Range("E36").Value = "  +5.69%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.920"
$c.ClearFormats()
$ws.Range("E37").Value = "  -5.40%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.900"
$c.ClearFormats()
$ws.Range("E38").Value = "  +5.78%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.52"
$c.ClearFormats()
$ws.Range("E39").Value = "  +6.53%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "37.29"
$c.ClearFormats()
$ws.Range("E40").Value = "  +1.88%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.65"
$c.ClearFormats()
$ws.Range("E41").Value = "  +0.07%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.644"
$c.ClearFormats()
$ws.Range("E42").Value = "  +6.79%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "20.26"
$c.ClearFormats()
$ws.Range("E43").Value = "  +2.34%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "278.10"
$c.ClearFormats()
$ws.Range("E44").Value = "  -1.30%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E45").Value = "  +0.21%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0987"
$c.ClearFormats()
$ws.Range("E46").Value = "  +0.12%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "4.95"
$c.ClearFormats()
$ws.Range("E47").Value = "  +5.60%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0544"
$c.ClearFormats()
$ws.Range("E48").Value = "  +3.28%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.093.83"
$c.ClearFormats()
$ws.Range("E49").Value = "  +0.35%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "10.55"
$c.ClearFormats()
$ws.Range("E50").Value = "  +2.34%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0233"
$c.ClearFormats()
$ws.Range("E51").Value = "  +1.18%  "
